$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 292.391276
$ws.Range("H2").Value = 877.173828
$ws.Range("I2").Value = 0.4546722242912879
$ws.Range("J2").Value = 0.4546722242912878
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5532856666666667
$ws.Range("N2").Value = 1.659857
$ws.Range("O2").Value = 0.4010144607159208
$ws.Range("P2").Value = 0.4010144607159208
$ws.Range("Q2").Value = 161.7759020691773
$ws.Range("R2").Value = 1455.983118622596
$ws.Range("S2").Value = 0.182330136826679
$ws.Range("T2").Value = 0.182330136826679
$ws.Range("G3").Value = 292.391276
$ws.Range("H3").Value = 877.173828
$ws.Range("I3").Value = 0.4546722242912879
$ws.Range("J3").Value = 0.4546722242912878
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8264293333333333
$ws.Range("N3").Value = 2.479288
$ws.Range("O3").Value = 0.5989855392840792
$ws.Range("P3").Value = 0.5989855392840792
$ws.Range("Q3").Value = 241.6407272971627
$ws.Range("R3").Value = 2174.766545674464
$ws.Range("S3").Value = 0.2723420874646089
$ws.Range("T3").Value = 0.2723420874646089
$ws.Range("H4").Value = 678.246018
$ws.Range("I4").Value = 0.3515604499097856
$ws.Range("J4").Value = 0.3515604499097856
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5532856666666667
$ws.Range("N4").Value = 1.659857
$ws.Range("O4").Value = 0.4010144607159208
$ws.Range("P4").Value = 0.4010144607159208
$ws.Range("Q4").Value = 125.0879334110473
$ws.Range("R4").Value = 1125.791400699426
$ws.Range("S4").Value = 0.1409808242296192
$ws.Range("T4").Value = 0.1409808242296192
$ws.Range("H5").Value = 678.246018
$ws.Range("I5").Value = 0.3515604499097856
$ws.Range("J5").Value = 0.3515604499097856
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8264293333333333
$ws.Range("N5").Value = 2.479288
$ws.Range("O5").Value = 0.5989855392840792
$ws.Range("P5").Value = 0.5989855392840792
$ws.Range("Q5").Value = 186.8408014972427
$ws.Range("R5").Value = 1681.567213475184
$ws.Range("S5").Value = 0.2105796256801664
$ws.Range("T5").Value = 0.2105796256801664
$ws.Range("G6").Value = 124.299764
$ws.Range("H6").Value = 372.899292
$ws.Range("I6").Value = 0.193287744251173
$ws.Range("J6").Value = 0.193287744251173
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5532856666666667
$ws.Range("N6").Value = 1.659857
$ws.Range("O6").Value = 0.4010144607159208
$ws.Range("P6").Value = 0.4010144607159208
$ws.Range("Q6").Value = 68.77327779124933
$ws.Range("R6").Value = 618.9595001212441
$ws.Range("S6").Value = 0.07751118052388098
$ws.Range("T6").Value = 0.07751118052388098
$ws.Range("G7").Value = 124.299764
$ws.Range("H7").Value = 372.899292
$ws.Range("I7").Value = 0.193287744251173
$ws.Range("J7").Value = 0.193287744251173
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8264293333333333
$ws.Range("N7").Value = 2.479288
$ws.Range("O7").Value = 0.5989855392840792
$ws.Range("P7").Value = 0.5989855392840792
$ws.Range("Q7").Value = 102.7249710960107
$ws.Range("R7").Value = 924.524739864096
$ws.Range("S7").Value = 0.1157765637272921
$ws.Range("T7").Value = 0.1157765637272921
$ws.Range("G8").Value = 0.30841
$ws.Range("H8").Value = 0.92523
$ws.Range("I8").Value = 0.0004795815477534156
$ws.Range("J8").Value = 0.0004795815477534155
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5532856666666667
$ws.Range("N8").Value = 1.659857
$ws.Range("O8").Value = 0.4010144607159208
$ws.Range("P8").Value = 0.4010144607159208
$ws.Range("Q8").Value = 0.1706388324566667
$ws.Range("R8").Value = 1.53574949211
$ws.Range("S8").Value = 0.0001923191357416425
$ws.Range("T8").Value = 0.0001923191357416425
$ws.Range("G9").Value = 0.30841
$ws.Range("H9").Value = 0.92523
$ws.Range("I9").Value = 0.0004795815477534156
$ws.Range("J9").Value = 0.0004795815477534155
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8264293333333333
$ws.Range("N9").Value = 2.479288
$ws.Range("O9").Value = 0.5989855392840792
$ws.Range("P9").Value = 0.5989855392840792
$ws.Range("Q9").Value = 0.2548790706933334
$ws.Range("R9").Value = 2.29391163624
$ws.Range("S9").Value = 0.000287262412011773
$ws.Range("T9").Value = 0.0002872624120117729
